# Vision bruger -> kunde
#
# 1. Change the stand-alone "Bruger" run (start of "Brugerprofil
#    håndtering:") to "Kunde".
# 2. The "_GoBack" bookmark (Word's "last edit location" bookmark) follows
#    the edit: remove it from its old spot at the very end of the document
#    and re-add it at the start of the "At oprette profil" paragraph,
#    mirroring where the edit above took place.

$d = $word.ActiveDocument

# --- 1. Bruger -> Kunde -----------------------------------------------
# Case-sensitive search so we hit only the capitalised stand-alone word
# (there are several lower-case "bruger" substrings elsewhere in the
# document that must stay untouched).
$rng = $d.Content
$found = $rng.Find.Execute("Bruger", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "Kunde"
}

# --- 2. Move the _GoBack bookmark --------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$target = $d.Content
$targetFound = $target.Find.Execute("At oprette profil")
if ($targetFound) {
    $insertionPoint = $d.Range($target.Start, $target.Start)
    $d.Bookmarks.Add("_GoBack", $insertionPoint)
}
